# SafetyChain test-data workbook refresh: replace the auto-generated
# document/task/form names with a new automation run's timestamps and
# add the QuestionaireForm submission column from the supplier portal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 existing columns: refresh auto-test names/timestamps ---
$ws.Range("A2").Value = "AUTO_TEST_FORM_ON_04/01/2019-18:43:51"
$ws.Range("B2").Value = "AUTO_TEST_TASK_ON_04/01/2019-18:50:31"
$ws.Range("C2").Value = "AUTO_TEST_DOCUMENT_ON_03/01/2019-14:11:31"
$ws.Range("D2").Value = "AUTO_TEST_RESTORE_DOCUMENT_ON_03/01/2019-14:11:31"
$ws.Range("E2").Value = "AUTO_TEST_ASSIGN_TASK_DOCUMENT_ON_03/01/2019-14:11:31"
$ws.Range("F2").Value = "AUTO_TEST_DMS_TASK_ON_03/01/2019-14:11:31"

# --- New column I: supplier-portal questionnaire submission ---
$ws.Range("I1").Value = "QuestionaireForm"
$ws.Range("I2").Value = "AUTO_TEST_FORM_ON_04/01/2019-08:59:33"
$ws.Range("I2").Font.Name = "Source Sans Pro"

# --- Column widths for the new/long text columns ---
$ws.Columns("F").ColumnWidth = 41.59
$ws.Columns("G").ColumnWidth = 30.1

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection ends on the newly added cell (drives dimension+selection) ---
$ws.Range("I2").Select() | Out-Null
